# Scheduled runner update: refresh currentAveragePrice / LevePrice / LeveProfit
# columns (H-N) across all class sheets with newly fetched market data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 231.54546
$ws.Range("I2").Value = 231.54546
$ws.Range("K2").Value = 231.54546
$ws.Range("M2").Value = -118.54546
$ws.Range("H41").Value = 610.3
$ws.Range("I41").Value = 474.75
$ws.Range("K41").Value = 474.75
$ws.Range("M41").Value = -34.75
$ws.Range("H64").Value = 63844.5
$ws.Range("J64").Value = 5990.8335
$ws.Range("L64").Value = 5990.8335
$ws.Range("N64").Value = -6486.8335
$ws.Range("H67").Value = 63844.5
$ws.Range("J67").Value = 5990.8335
$ws.Range("L67").Value = 5990.8335
$ws.Range("N67").Value = -7706.8335
$ws.Range("H98").Value = 27077.26
$ws.Range("I98").Value = 38611.066
$ws.Range("K98").Value = 38611.066
$ws.Range("M98").Value = -37113.066
$ws.Range("H100").Value = 65432.086
$ws.Range("I100").Value = 72154.14
$ws.Range("J100").Value = 54975.555
$ws.Range("K100").Value = 72154.14
$ws.Range("L100").Value = 54975.555
$ws.Range("M100").Value = -71613.14
$ws.Range("N100").Value = -56057.555
$ws.Range("H122").Value = 27077.26
$ws.Range("I122").Value = 38611.066
$ws.Range("K122").Value = 115833.198
$ws.Range("M122").Value = -113383.198
$ws.Range("H132").Value = 1926070.6
$ws.Range("I132").Value = 2840.4255
$ws.Range("K132").Value = 8521.2765
$ws.Range("M132").Value = -5991.2765

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7153.887
$ws.Range("I32").Value = 7043.044
$ws.Range("J32").Value = 9666.333000000001
$ws.Range("K32").Value = 7043.044
$ws.Range("L32").Value = 9666.333000000001
$ws.Range("M32").Value = -6756.044
$ws.Range("N32").Value = -10240.333
$ws.Range("H61").Value = 10663.429
$ws.Range("I61").Value = 11607.556
$ws.Range("J61").Value = 4998.6665
$ws.Range("K61").Value = 11607.556
$ws.Range("L61").Value = 4998.6665
$ws.Range("M61").Value = -11395.556
$ws.Range("N61").Value = -5422.6665
$ws.Range("H74").Value = 6162.3335
$ws.Range("I74").Value = 12077.333
$ws.Range("K74").Value = 12077.333
$ws.Range("M74").Value = -11203.333
$ws.Range("H77").Value = 6162.3335
$ws.Range("I77").Value = 12077.333
$ws.Range("K77").Value = 60386.665
$ws.Range("M77").Value = -56018.665
$ws.Range("H88").Value = 1562.1538
$ws.Range("I88").Value = 932.3333
$ws.Range("K88").Value = 932.3333
$ws.Range("M88").Value = -526.3333
$ws.Range("H91").Value = 1562.1538
$ws.Range("I91").Value = 932.3333
$ws.Range("K91").Value = 932.3333
$ws.Range("M91").Value = 471.6667
$ws.Range("H122").Value = 1370061.5
$ws.Range("I122").Value = 6186.3335
$ws.Range("K122").Value = 18559.0005
$ws.Range("M122").Value = -16109.0005
$ws.Range("H132").Value = 3500.4546
$ws.Range("I132").Value = 2501.1428
$ws.Range("K132").Value = 7503.428400000001
$ws.Range("M132").Value = -4973.428400000001
$ws.Range("H136").Value = 10663.429
$ws.Range("I136").Value = 11607.556
$ws.Range("J136").Value = 4998.6665
$ws.Range("K136").Value = 34822.66800000001
$ws.Range("L136").Value = 14995.9995
$ws.Range("M136").Value = -32272.66800000001
$ws.Range("N136").Value = -20095.9995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 7893.1333
$ws.Range("I86").Value = 12024.625
$ws.Range("J86").Value = 3171.4285
$ws.Range("K86").Value = 12024.625
$ws.Range("L86").Value = 3171.4285
$ws.Range("M86").Value = -10901.625
$ws.Range("N86").Value = -5417.4285
$ws.Range("H89").Value = 7893.1333
$ws.Range("I89").Value = 12024.625
$ws.Range("J89").Value = 3171.4285
$ws.Range("K89").Value = 60123.125
$ws.Range("L89").Value = 15857.1425
$ws.Range("M89").Value = -54507.125
$ws.Range("N89").Value = -27089.1425
$ws.Range("H94").Value = 2500.487
$ws.Range("J94").Value = 5271.8
$ws.Range("L94").Value = 5271.8
$ws.Range("N94").Value = -6173.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6134.8184
$ws.Range("I31").Value = 6090.7407
$ws.Range("K31").Value = 6090.7407
$ws.Range("M31").Value = -5795.7407
$ws.Range("H34").Value = 6134.8184
$ws.Range("I34").Value = 6090.7407
$ws.Range("K34").Value = 6090.7407
$ws.Range("M34").Value = -5888.7407
$ws.Range("H47").Value = 46935
$ws.Range("J47").Value = 46935
$ws.Range("L47").Value = 46935
$ws.Range("N47").Value = -48067
$ws.Range("H58").Value = 5629.4546
$ws.Range("I58").Value = 6858.7144
$ws.Range("K58").Value = 6858.7144
$ws.Range("M58").Value = -6655.7144
$ws.Range("H122").Value = 6644.7144
$ws.Range("I122").Value = 9281.5
$ws.Range("J122").Value = 1371.1428
$ws.Range("K122").Value = 27844.5
$ws.Range("L122").Value = 4113.428400000001
$ws.Range("M122").Value = -25394.5
$ws.Range("N122").Value = -9013.428400000001
$ws.Range("H132").Value = 1825.5264
$ws.Range("I132").Value = 1391.1538
$ws.Range("J132").Value = 2766.6667
$ws.Range("K132").Value = 4173.4614
$ws.Range("L132").Value = 8300.000100000001
$ws.Range("M132").Value = -1643.4614
$ws.Range("N132").Value = -13360.0001
$ws.Range("H134").Value = 12875.272
$ws.Range("I134").Value = 15191.333
$ws.Range("K134").Value = 45573.999
$ws.Range("M134").Value = -43038.999
$ws.Range("H136").Value = 5629.4546
$ws.Range("I136").Value = 6858.7144
$ws.Range("K136").Value = 20576.1432
$ws.Range("M136").Value = -18026.1432

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("M31").ClearContents()
$ws.Range("H44").Value = 2835.8333
$ws.Range("I44").Value = 403
$ws.Range("J44").Value = 15000
$ws.Range("K44").Value = 1209
$ws.Range("L44").Value = 45000
$ws.Range("M44").Value = -811
$ws.Range("N44").Value = -45796
$ws.Range("H103").Value = 3419.8
$ws.Range("I103").Value = 4644.75
$ws.Range("K103").Value = 13934.25
$ws.Range("M103").Value = -13055.25
$ws.Range("H119").Value = 776
$ws.Range("I119").Value = 776
$ws.Range("K119").Value = 2328
$ws.Range("M119").Value = 2510
$ws.Range("H121").Value = 5116.6665
$ws.Range("I121").Value = 9000
$ws.Range("J121").Value = 3175
$ws.Range("K121").Value = 27000
$ws.Range("L121").Value = 9525
$ws.Range("M121").Value = -25690
$ws.Range("N121").Value = -12145
$ws.Range("H134").Value = 2893.5833
$ws.Range("I134").Value = 2272.3
$ws.Range("K134").Value = 6816.900000000001
$ws.Range("M134").Value = -1746.900000000001
$ws.Range("H139").Value = 3335839.8
$ws.Range("I139").Value = 5001759.5
$ws.Range("K139").Value = 15005278.5
$ws.Range("M139").Value = -15000138.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 5007.6562
$ws.Range("I102").Value = 4958.1665
$ws.Range("K102").Value = 4958.1665
$ws.Range("M102").Value = -3336.1665
$ws.Range("H132").Value = 4264.9355
$ws.Range("I132").Value = 4467
$ws.Range("J132").Value = 2379
$ws.Range("K132").Value = 13401
$ws.Range("L132").Value = 7137
$ws.Range("M132").Value = -10871
$ws.Range("N132").Value = -12197

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 61935.125
$ws.Range("I40").Value = 72747.836
$ws.Range("K40").Value = 72747.836
$ws.Range("M40").Value = -72611.836
$ws.Range("H46").Value = 4514.9
$ws.Range("J46").Value = 6999.8335
$ws.Range("L46").Value = 6999.8335
$ws.Range("N46").Value = -7375.8335
$ws.Range("H122").Value = 5656.3887
$ws.Range("I122").Value = 6012
$ws.Range("K122").Value = 18036
$ws.Range("M122").Value = -15586

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1206.3422
$ws.Range("I113").Value = 677.6429000000001
$ws.Range("K113").Value = 2032.9287
$ws.Range("M113").Value = 137.0712999999998
$ws.Range("H132").Value = 4031.827
$ws.Range("I132").Value = 3868.889
$ws.Range("K132").Value = 11606.667
$ws.Range("M132").Value = -9076.667000000001
$ws.Range("H136").Value = 347506.34
$ws.Range("J136").Value = 11983.333
$ws.Range("L136").Value = 35949.999
$ws.Range("N136").Value = -41049.999
